$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.756.16"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "2.333.79"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  -3.37%  "
$c = $ws.Range("D5")
$c.Value = "'313.36"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "
$c = $ws.Range("D6")
$c.Value = "'107.88"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.52%  "
$c = $ws.Range("D7")
$c.Value = "'0.631"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("E8").Value = "  -0.32%  "
$c = $ws.Range("D9")
$c.Value = "'0.620"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.90%  "
$c = $ws.Range("D10")
$c.Value = "'41.32"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +6.42%  "
$c = $ws.Range("D11")
$c.Value = "'0.0918"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("E12").Value = "  +3.92%  "
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("E14").Value = "  -0.42%  "
$c = $ws.Range("D15")
$c.Value = "'15.47"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("D16").Value = "2.686.36"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("D17").Value = "2.328.93"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").Value = "43.689.08"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("E20").Value = "  +2.55%  "
$c = $ws.Range("D21")
$c.Value = "'13.00"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.22%  "
$c = $ws.Range("D22")
$c.Value = "'74.34"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.41%  "
$c = $ws.Range("D23")
$c.Value = "'3.48"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.77%  "
$c = $ws.Range("D24")
$c.Value = "'268.70"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("E25").Value = "  +4.84%  "
$ws.Range("E26").Value = "  -0.18%  "
$c = $ws.Range("D27")
$c.Value = "'7.62"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +10.78%  "
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("E29").Value = "  -1.38%  "
$c = $ws.Range("D30")
$c.Value = "'39.53"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +9.42%  "
$c = $ws.Range("D31")
$c.Value = "'22.54"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.74%  "
$c = $ws.Range("D32")
$c.Value = "'168.09"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.95%  "
$c = $ws.Range("D33")
$c.Value = "'0.0900"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.68%  "
$c = $ws.Range("D34")
$c.Value = "'2.86"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +8.85%  "
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("E36").Value = "  +3.46%  "
$c = $ws.Range("D37")
$c.Value = "'4.70"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.38%  "
$ws.Range("E38").Value = "  +5.18%  "
$c = $ws.Range("D39")
$c.Value = "'2.90"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +10.13%  "
$c = $ws.Range("D40")
$c.Value = "'3.78"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.12%  "
$ws.Range("E41").Value = "  +10.12%  "
$c = $ws.Range("D42")
$c.Value = "'104.03"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +12.07%  "
$c = $ws.Range("D43")
$c.Value = "'13.78"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +15.83%  "
$ws.Range("E44").Value = "  +5.93%  "
$c = $ws.Range("D45")
$c.Value = "'71.71"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.55%  "
$ws.Range("E46").Value = "  -0.13%  "
$c = $ws.Range("D47")
$c.Value = "'115.03"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.52%  "
$ws.Range("E48").Value = "  +17.85%  "
$ws.Range("D49").Value = "1.656.94"
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("E50").Value = "  +3.81%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$c = $ws.Range("D51")
$c.Value = "'75.95"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.74%  "
